$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The 2013 section had a stray blank row (old row 35) between its header
# ("2013") and the data table that starts with the FY-column header row.
# Delete that entire row so everything below (the 2013/2014/2015 schedule
# tables) shifts up by one row, matching the rest of the sheet's layout
# (each year's header row is immediately followed by its column-header row).
$ws.Rows("35").Delete()

# Leave the selection on the row that is now in the position of the old
# row 35 (the FY column-header row for the 2013 table), matching where the
# cursor ended up after the delete.
$ws.Range("A35:XFD35").Select() | Out-Null
